$wb = $excel.ActiveWorkbook

# Update the quiz-name content on the "Admin_messaging" sheet:
# the placeholder student name "virat" becomes "ronaldo" (and the
# dependent message text updates to match).
$wsAdminMessaging = $wb.Worksheets.Item("Admin_messaging")
$wsAdminMessaging.Range("C2").Value = "ronaldo"
$wsAdminMessaging.Range("D2").Value = "hi ronaldo how are you"

# Move the active selection on that sheet from D7 to D6.
$wsAdminMessaging.Range("D6").Select()

# "Admin_messaging" becomes the active/selected sheet (it also becomes the
# active tab stored on the workbook), taking over from "deleting_student".
$wsAdminMessaging.Activate()
